$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "BK308 / Odborný seminář B / Kadlečková Ivana" record (old row 2) was
# dropped and the rows below it moved up by one. A couple of other values
# were also corrected independently of that shift:
#  - row 3 (KRSP -> KKRSP / KRSPKBI -> KKRSPKBI)
#  - the date on what is now row 6 (2.11.2024 -> 16.11.2024)
#  - the surrounding quote marks around the "seminarici" teacher names were
#    stripped everywhere.
#
# Force the date-like and purely-numeric-looking text columns to stay text
# (so Excel doesn't reinterpret "6.12.2024" as a date serial or "1543" as a
# number) before writing the values into them.
$ws.Range("F2:G4").NumberFormat = "@"
$ws.Range("F6:G6").NumberFormat = "@"
$ws.Range("N2:N3").NumberFormat = "@"

# Row 2 (was row 3): BK320 / Seminář z průmyslové chemie
$ws.Cells.Item(2, 2).Value = "BK320"
$ws.Cells.Item(2, 3).Value = 9475
$ws.Cells.Item(2, 6).Value = "6.12.2024"
$ws.Cells.Item(2, 7).Value = "6.12.2024"
$ws.Cells.Item(2, 8).Value = "10:00"
$ws.Cells.Item(2, 9).Value = "13:50"
$ws.Cells.Item(2, 10).Value = "BK320KCH"
$ws.Cells.Item(2, 12).Value = "Seminář z průmyslové chemie"
$ws.Cells.Item(2, 13).Value = "doc. Ing. Jaromír Lederer, CSc."
$ws.Cells.Item(2, 14).Value = "1543"
$ws.Cells.Item(2, 15).Value = "Šimek Josef, Ing. Ph.D."

# Row 3 (was row 4): KBI / KKRSP (note: zkratka/identifier differ from the
# unchanged KRSP rows below)
$ws.Cells.Item(3, 1).Value = "KBI"
$ws.Cells.Item(3, 2).Value = "KKRSP"
$ws.Cells.Item(3, 3).Value = 8952
$ws.Cells.Item(3, 6).Value = "16.11.2024"
$ws.Cells.Item(3, 7).Value = "16.11.2024"
$ws.Cells.Item(3, 8).Value = "09:00"
$ws.Cells.Item(3, 9).Value = "12:50"
$ws.Cells.Item(3, 10).Value = "KKRSPKBI"
$ws.Cells.Item(3, 11).Value = "KBI"
$ws.Cells.Item(3, 12).Value = "Komplexní reflektivní seminář praxe"
$ws.Cells.Item(3, 13).Value = "doc. PhDr. Kateřina Jančaříková, Ph.D."
$ws.Cells.Item(3, 14).Value = "8441"
$ws.Cells.Item(3, 15).Value = "Svobodová Silvie, PhDr. Ing. Ph.D."

# Row 4 (was row 5): date shifted, quotes stripped
$ws.Cells.Item(4, 6).Value = "2.11.2024"
$ws.Cells.Item(4, 7).Value = "2.11.2024"
$ws.Cells.Item(4, 13).Value = "doc. PhDr. Kateřina Jančaříková, Ph.D."

# Row 5 (was row 6, unchanged dates): quotes stripped
$ws.Cells.Item(5, 13).Value = "doc. PhDr. Kateřina Jančaříková, Ph.D."

# Row 6 (was row 6 originally, but date corrected, quotes stripped)
$ws.Cells.Item(6, 6).Value = "16.11.2024"
$ws.Cells.Item(6, 7).Value = "16.11.2024"
$ws.Cells.Item(6, 13).Value = "doc. PhDr. Kateřina Jančaříková, Ph.D."

# Row 7: quotes stripped
$ws.Cells.Item(7, 13).Value = "RNDr. Jiří Králík, Ph.D."
